$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 127.181816
$ws.Range("I4").Value = 89.90000000000001
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 89.90000000000001
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 24.09999999999999
$ws.Range("N4").Value = -728
$ws.Range("H41").Value = 1222.9231
$ws.Range("I41").Value = 222.5
$ws.Range("J41").Value = 1667.5555
$ws.Range("K41").Value = 222.5
$ws.Range("L41").Value = 1667.5555
$ws.Range("M41").Value = 217.5
$ws.Range("N41").Value = -2547.5555
$ws.Range("H106").Value = 2778.1
$ws.Range("I106").Value = 2120.111
$ws.Range("K106").Value = 2120.111
$ws.Range("M106").Value = -1489.111
$ws.Range("H111").Value = 1496.25
$ws.Range("I111").Value = 1277.8
$ws.Range("J111").Value = 1860.3334
$ws.Range("K111").Value = 3833.4
$ws.Range("L111").Value = 5581.0002
$ws.Range("M111").Value = -766.3999999999996
$ws.Range("N111").Value = -11715.0002
$ws.Range("H112").Value = 10189.929
$ws.Range("I112").Value = 966.3333
$ws.Range("K112").Value = 2898.9999
$ws.Range("M112").Value = -1790.9999
$ws.Range("H113").Value = 2728.2
$ws.Range("I113").Value = 2728.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2728.2
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 525.8000000000002
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 3402837.5
$ws.Range("I132").Value = 4082908.8
$ws.Range("K132").Value = 12248726.4
$ws.Range("M132").Value = -12246196.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 573.8
$ws.Range("I2").Value = 493
$ws.Range("J2").Value = 695
$ws.Range("K2").Value = 493
$ws.Range("L2").Value = 695
$ws.Range("M2").Value = -380
$ws.Range("N2").Value = -921
$ws.Range("H61").Value = 7927.3447
$ws.Range("I61").Value = 4786.696
$ws.Range("K61").Value = 4786.696
$ws.Range("M61").Value = -4574.696
$ws.Range("H110").Value = 20835.38
$ws.Range("J110").Value = 2750.6
$ws.Range("L110").Value = 2750.6
$ws.Range("N110").Value = -6840.6
$ws.Range("H116").Value = 573.8
$ws.Range("I116").Value = 493
$ws.Range("J116").Value = 695
$ws.Range("K116").Value = 493
$ws.Range("L116").Value = 695
$ws.Range("M116").Value = 1801
$ws.Range("N116").Value = -5283
$ws.Range("H128").Value = 166666.67
$ws.Range("J128").Value = 166666.67
$ws.Range("L128").Value = 166666.67
$ws.Range("N128").Value = -176626.67
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 1986.2424
$ws.Range("I132").Value = 1594.037
$ws.Range("K132").Value = 4782.111
$ws.Range("M132").Value = -2252.111
$ws.Range("H135").Value = 91332
$ws.Range("J135").Value = 91332
$ws.Range("L135").Value = 91332
$ws.Range("N135").Value = -101472
$ws.Range("H136").Value = 7927.3447
$ws.Range("I136").Value = 4786.696
$ws.Range("K136").Value = 14360.088
$ws.Range("M136").Value = -11810.088
$ws.Range("H139").Value = 116374.75
$ws.Range("J139").Value = 116374.75
$ws.Range("L139").Value = 116374.75
$ws.Range("N139").Value = -126654.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 573.8
$ws.Range("I3").Value = 493
$ws.Range("J3").Value = 695
$ws.Range("K3").Value = 493
$ws.Range("L3").Value = 695
$ws.Range("M3").Value = -379
$ws.Range("N3").Value = -923
$ws.Range("H74").Value = 90000
$ws.Range("J74").Value = 90000
$ws.Range("L74").Value = 90000
$ws.Range("N74").Value = -91872
$ws.Range("H77").Value = 90000
$ws.Range("J77").Value = 90000
$ws.Range("L77").Value = 270000
$ws.Range("N77").Value = -279360
$ws.Range("H134").Value = 3843.8
$ws.Range("I134").Value = 3781.12
$ws.Range("J134").Value = 4157.2
$ws.Range("K134").Value = 11343.36
$ws.Range("L134").Value = 12471.6
$ws.Range("M134").Value = -8808.360000000001
$ws.Range("N134").Value = -17541.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8336666.5
$ws.Range("J4").Value = 8336666.5
$ws.Range("L4").Value = 8336666.5
$ws.Range("N4").Value = -8336890.5
$ws.Range("H58").Value = 1581.3334
$ws.Range("J58").Value = 2032.6666
$ws.Range("L58").Value = 2032.6666
$ws.Range("N58").Value = -2438.6666
$ws.Range("H132").Value = 2927.5881
$ws.Range("I132").Value = 2125.818
$ws.Range("J132").Value = 4397.5
$ws.Range("K132").Value = 6377.454000000001
$ws.Range("L132").Value = 13192.5
$ws.Range("M132").Value = -3847.454000000001
$ws.Range("N132").Value = -18252.5
$ws.Range("H134").Value = 2243.6
$ws.Range("I134").Value = 2243.6
$ws.Range("K134").Value = 6730.799999999999
$ws.Range("M134").Value = -4195.799999999999
$ws.Range("H136").Value = 1581.3334
$ws.Range("J136").Value = 2032.6666
$ws.Range("L136").Value = 6097.9998
$ws.Range("N136").Value = -11197.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6120.9165
$ws.Range("I56").Value = 6120.9165
$ws.Range("K56").Value = 6120.9165
$ws.Range("M56").Value = -5590.9165
$ws.Range("H94").Value = 5155.1113
$ws.Range("I94").Value = 300
$ws.Range("K94").Value = 900
$ws.Range("M94").Value = -224
$ws.Range("H115").Value = 2020
$ws.Range("J115").Value = 2020
$ws.Range("L115").Value = 6060
$ws.Range("N115").Value = -8410
$ws.Range("H131").Value = 2553.4773
$ws.Range("J131").Value = 2049.818
$ws.Range("L131").Value = 6149.454000000001
$ws.Range("N131").Value = -16229.454
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 33152.375
$ws.Range("J20").Value = 33152.375
$ws.Range("L20").Value = 33152.375
$ws.Range("N20").Value = -33642.375
$ws.Range("H24").Value = 28490.158
$ws.Range("J24").Value = 28490.158
$ws.Range("L24").Value = 28490.158
$ws.Range("N24").Value = -28836.158
$ws.Range("H122").Value = 33335122
$ws.Range("I122").Value = 1011.75
$ws.Range("J122").Value = 71431250
$ws.Range("K122").Value = 3035.25
$ws.Range("L122").Value = 214293750
$ws.Range("M122").Value = -585.25
$ws.Range("N122").Value = -214298650
$ws.Range("H126").Value = 1779.9
$ws.Range("I126").Value = 960
$ws.Range("K126").Value = 2880
$ws.Range("M126").Value = -410
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1599.7142
$ws.Range("I16").Value = 1450
$ws.Range("J16").Value = 1799.3334
$ws.Range("K16").Value = 1450
$ws.Range("L16").Value = 1799.3334
$ws.Range("M16").Value = -1280
$ws.Range("N16").Value = -2139.3334
$ws.Range("H46").Value = 3671.7083
$ws.Range("I46").Value = 1667.6666
$ws.Range("J46").Value = 3958
$ws.Range("K46").Value = 1667.6666
$ws.Range("L46").Value = 3958
$ws.Range("M46").Value = -1479.6666
$ws.Range("N46").Value = -4334
$ws.Range("H63").Value = 85711.86
$ws.Range("J63").Value = 85711.86
$ws.Range("L63").Value = 85711.86
$ws.Range("N63").Value = -87209.86
$ws.Range("H66").Value = 85711.86
$ws.Range("J66").Value = 85711.86
$ws.Range("L66").Value = 257135.58
$ws.Range("N66").Value = -264623.58
$ws.Range("H86").Value = 112500
$ws.Range("J86").Value = 112500
$ws.Range("L86").Value = 112500
$ws.Range("N86").Value = -114872
$ws.Range("H89").Value = 112500
$ws.Range("J89").Value = 112500
$ws.Range("L89").Value = 337500
$ws.Range("N89").Value = -349356
$ws.Range("H122").Value = 4395.8
$ws.Range("I122").Value = 3994.6667
$ws.Range("K122").Value = 11984.0001
$ws.Range("M122").Value = -9534.000100000001
$ws.Range("H132").Value = 5154.037
$ws.Range("I132").Value = 5511.231
$ws.Range("J132").Value = 4822.357
$ws.Range("K132").Value = 16533.693
$ws.Range("L132").Value = 14467.071
$ws.Range("M132").Value = -14003.693
$ws.Range("N132").Value = -19527.071
$ws.Range("H133").Value = 109999
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 142859470
$ws.Range("I96").Value = 1400
$ws.Range("J96").Value = 200002700
$ws.Range("K96").Value = 1400
$ws.Range("L96").Value = 200002700
$ws.Range("M96").Value = -27
$ws.Range("N96").Value = -200005446
$ws.Range("H132").Value = 9287072
$ws.Range("I132").Value = 10447510
$ws.Range("K132").Value = 31342530
$ws.Range("M132").Value = -31340000
$ws.Range("H136").Value = 7413.8867
$ws.Range("I136").Value = 9122.424000000001
$ws.Range("J136").Value = 4594.8
$ws.Range("K136").Value = 27367.272
$ws.Range("L136").Value = 13784.4
$ws.Range("M136").Value = -24817.272
$ws.Range("N136").Value = -18884.4
